$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4998.3335
$ws.Range("I32").Value = 4997.5
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 4997.5
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -4671.5
$ws.Range("N32").Value = -5652
$ws.Range("H34").Value = 4108.143
$ws.Range("I34").Value = 4108.143
$ws.Range("K34").Value = 4108.143
$ws.Range("M34").Value = -3905.143
$ws.Range("H36").Value = 4108.143
$ws.Range("I36").Value = 4108.143
$ws.Range("K36").Value = 4108.143
$ws.Range("M36").Value = -3393.143
$ws.Range("H40").Value = 2170.6428
$ws.Range("I40").Value = 1914.2858
$ws.Range("J40").Value = 2427
$ws.Range("K40").Value = 1914.2858
$ws.Range("L40").Value = 2427
$ws.Range("M40").Value = -1739.2858
$ws.Range("N40").Value = -2777
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = $null
$ws.Range("H62").Value = 7333
$ws.Range("I62").Value = 6799.6
$ws.Range("K62").Value = 6799.6
$ws.Range("M62").Value = -6175.6
$ws.Range("H65").Value = 7333
$ws.Range("I65").Value = 6799.6
$ws.Range("K65").Value = 33998
$ws.Range("M65").Value = -30878
$ws.Range("H69").Value = 7500
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14126
$ws.Range("H72").Value = 7500
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -40632
$ws.Range("H111").Value = 2025.7727
$ws.Range("I111").Value = 1056.0667
$ws.Range("J111").Value = 4103.7144
$ws.Range("K111").Value = 3168.2001
$ws.Range("L111").Value = 12311.1432
$ws.Range("M111").Value = -101.2001
$ws.Range("N111").Value = -18445.1432
$ws.Range("H132").Value = 2284.2727
$ws.Range("I132").Value = 2055.5789
$ws.Range("J132").Value = 3732.6667
$ws.Range("K132").Value = 6166.736699999999
$ws.Range("L132").Value = 11198.0001
$ws.Range("M132").Value = -3636.736699999999
$ws.Range("N132").Value = -16258.0001
$ws.Range("H135").Value = 1335
$ws.Range("I135").Value = 1114.5454
$ws.Range("K135").Value = 10030.9086
$ws.Range("M135").Value = -7495.908599999999
$ws.Range("H137").Value = 2121.4666
$ws.Range("I137").Value = 2126.8333
$ws.Range("J137").Value = 2100
$ws.Range("K137").Value = 6380.499899999999
$ws.Range("L137").Value = 6300
$ws.Range("M137").Value = -3830.499899999999
$ws.Range("N137").Value = -11400
$ws.Range("H141").Value = 5102.6665
$ws.Range("I141").Value = 3373.2
$ws.Range("J141").Value = 13750
$ws.Range("K141").Value = 10119.6
$ws.Range("L141").Value = 41250
$ws.Range("M141").Value = -4939.599999999999
$ws.Range("N141").Value = -51610

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2495.3333
$ws.Range("I2").Value = 1621.3636
$ws.Range("K2").Value = 1621.3636
$ws.Range("M2").Value = -1508.3636
$ws.Range("H27").Value = 5000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = -5368
$ws.Range("H45").Value = 1466.6666
$ws.Range("I45").Value = 1466.6666
$ws.Range("K45").Value = 1466.6666
$ws.Range("M45").Value = -1089.6666
$ws.Range("H46").Value = 15481.714
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null
$ws.Range("H50").Value = 400
$ws.Range("I50").Value = 523.5
$ws.Range("J50").Value = 276.5
$ws.Range("K50").Value = 523.5
$ws.Range("L50").Value = 276.5
$ws.Range("M50").Value = 190.5
$ws.Range("N50").Value = -1704.5
$ws.Range("H61").Value = 2899.875
$ws.Range("I61").Value = 2955.1428
$ws.Range("J61").Value = 2513
$ws.Range("K61").Value = 2955.1428
$ws.Range("L61").Value = 2513
$ws.Range("M61").Value = -2743.1428
$ws.Range("N61").Value = -2937
$ws.Range("H74").Value = 1399.862
$ws.Range("I74").Value = 1201.6428
$ws.Range("K74").Value = 1201.6428
$ws.Range("M74").Value = -327.6428000000001
$ws.Range("H77").Value = 1399.862
$ws.Range("I77").Value = 1201.6428
$ws.Range("K77").Value = 6008.214
$ws.Range("M77").Value = -1640.214
$ws.Range("H88").Value = 512.41174
$ws.Range("J88").Value = 760.2222
$ws.Range("L88").Value = 760.2222
$ws.Range("N88").Value = -1572.2222
$ws.Range("H91").Value = 512.41174
$ws.Range("J91").Value = 760.2222
$ws.Range("L91").Value = 760.2222
$ws.Range("N91").Value = -3568.2222
$ws.Range("H95").Value = 45098.57
$ws.Range("J95").Value = 45098.57
$ws.Range("L95").Value = 45098.57
$ws.Range("N95").Value = -50590.57
$ws.Range("H102").Value = 600.375
$ws.Range("I102").Value = 634
$ws.Range("K102").Value = 634
$ws.Range("M102").Value = 988
$ws.Range("H116").Value = 2495.3333
$ws.Range("I116").Value = 1621.3636
$ws.Range("K116").Value = 1621.3636
$ws.Range("M116").Value = 672.6364000000001
$ws.Range("H122").Value = 457505.38
$ws.Range("I122").Value = 669554.6
$ws.Range("K122").Value = 2008663.8
$ws.Range("M122").Value = -2006213.8
$ws.Range("H132").Value = 1410.5
$ws.Range("I132").Value = 1410.5
$ws.Range("K132").Value = 4231.5
$ws.Range("M132").Value = -1701.5
$ws.Range("H136").Value = 2899.875
$ws.Range("I136").Value = 2955.1428
$ws.Range("J136").Value = 2513
$ws.Range("K136").Value = 8865.428400000001
$ws.Range("L136").Value = 7539
$ws.Range("M136").Value = -6315.428400000001
$ws.Range("N136").Value = -12639

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2495.3333
$ws.Range("I3").Value = 1621.3636
$ws.Range("K3").Value = 1621.3636
$ws.Range("M3").Value = -1507.3636
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184
$ws.Range("H86").Value = 5801
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 6701.5
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 6701.5
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -8947.5
$ws.Range("H89").Value = 5801
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 6701.5
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 33507.5
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -44739.5
$ws.Range("H94").Value = 645.65216
$ws.Range("I94").Value = 676.5
$ws.Range("K94").Value = 676.5
$ws.Range("M94").Value = -225.5
$ws.Range("H99").Value = 31760.242
$ws.Range("I99").Value = 39802.96
$ws.Range("K99").Value = 39802.96
$ws.Range("M99").Value = -38304.96
$ws.Range("H105").Value = 1000.0455
$ws.Range("I105").Value = 880.05
$ws.Range("K105").Value = 880.05
$ws.Range("M105").Value = 866.95
$ws.Range("H107").Value = 1013.3333
$ws.Range("I107").Value = 920
$ws.Range("K107").Value = 920
$ws.Range("M107").Value = 1000
$ws.Range("H130").Value = 62446.668
$ws.Range("J130").Value = 62446.668
$ws.Range("L130").Value = 62446.668
$ws.Range("N130").Value = -72486.66800000001
$ws.Range("H134").Value = 2422.3333
$ws.Range("I134").Value = 2107
$ws.Range("K134").Value = 6321
$ws.Range("M134").Value = -3786

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 574.5
$ws.Range("I16").Value = 574.5
$ws.Range("K16").Value = 574.5
$ws.Range("M16").Value = -287.5
$ws.Range("H31").Value = 3100.8044
$ws.Range("I31").Value = 2793.04
$ws.Range("J31").Value = 3467.1904
$ws.Range("K31").Value = 2793.04
$ws.Range("L31").Value = 3467.1904
$ws.Range("M31").Value = -2498.04
$ws.Range("N31").Value = -4057.1904
$ws.Range("H34").Value = 3100.8044
$ws.Range("I34").Value = 2793.04
$ws.Range("J34").Value = 3467.1904
$ws.Range("K34").Value = 2793.04
$ws.Range("L34").Value = 3467.1904
$ws.Range("M34").Value = -2591.04
$ws.Range("N34").Value = -3871.1904
$ws.Range("H86").Value = 9786.1
$ws.Range("I86").Value = 8295.200000000001
$ws.Range("J86").Value = 11277
$ws.Range("K86").Value = 8295.200000000001
$ws.Range("L86").Value = 11277
$ws.Range("M86").Value = -7172.200000000001
$ws.Range("N86").Value = -13523
$ws.Range("H89").Value = 9786.1
$ws.Range("I89").Value = 8295.200000000001
$ws.Range("J89").Value = 11277
$ws.Range("K89").Value = 41476
$ws.Range("L89").Value = 56385
$ws.Range("M89").Value = -35860
$ws.Range("N89").Value = -67617
$ws.Range("H105").Value = 1223.4
$ws.Range("I105").Value = 699
$ws.Range("J105").Value = 1573
$ws.Range("K105").Value = 699
$ws.Range("L105").Value = 1573
$ws.Range("M105").Value = 1048
$ws.Range("N105").Value = -5067
$ws.Range("H107").Value = 1072.25
$ws.Range("I107").Value = 820.1667
$ws.Range("K107").Value = 820.1667
$ws.Range("M107").Value = 1099.8333
$ws.Range("H113").Value = 574.5
$ws.Range("I113").Value = 574.5
$ws.Range("K113").Value = 574.5
$ws.Range("M113").Value = 1595.5
$ws.Range("H132").Value = 2052.4285
$ws.Range("I132").Value = 1680.9231
$ws.Range("J132").Value = 6882
$ws.Range("K132").Value = 5042.7693
$ws.Range("L132").Value = 20646
$ws.Range("M132").Value = -2512.7693
$ws.Range("N132").Value = -25706
$ws.Range("H134").Value = 1683.579
$ws.Range("I134").Value = 1396.721
$ws.Range("K134").Value = 4190.163
$ws.Range("M134").Value = -1655.163
$ws.Range("H135").Value = 130000
$ws.Range("J135").Value = 130000
$ws.Range("L135").Value = 130000
$ws.Range("N135").Value = -140140

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50026.7
$ws.Range("I2").Value = 83359.164
$ws.Range("J2").Value = 28
$ws.Range("K2").Value = 500154.9840000001
$ws.Range("L2").Value = 168
$ws.Range("M2").Value = -500041.9840000001
$ws.Range("N2").Value = -394
$ws.Range("H11").Value = 150583.75
$ws.Range("I11").Value = 150583.75
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 451751.25
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -451611.25
$ws.Range("N11").Value = $null
$ws.Range("H38").Value = 132.42857
$ws.Range("I38").Value = 85.75
$ws.Range("K38").Value = 257.25
$ws.Range("M38").Value = 89.75
$ws.Range("H68").Value = 350.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("H71").Value = 350.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
$ws.Range("H92").Value = 420.375
$ws.Range("I92").Value = 442.41666
$ws.Range("K92").Value = 1327.24998
$ws.Range("M92").Value = -79.24998000000005
$ws.Range("H113").Value = 2021.125
$ws.Range("I113").Value = 2190
$ws.Range("J113").Value = 1919.8
$ws.Range("K113").Value = 6570
$ws.Range("L113").Value = 5759.4
$ws.Range("M113").Value = -4400
$ws.Range("N113").Value = -10099.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7785
$ws.Range("I70").Value = 7498.3335
$ws.Range("K70").Value = 7498.3335
$ws.Range("M70").Value = -7228.3335
$ws.Range("H73").Value = 7785
$ws.Range("I73").Value = 7498.3335
$ws.Range("K73").Value = 7498.3335
$ws.Range("M73").Value = -6562.3335
$ws.Range("H80").Value = 2012.421
$ws.Range("I80").Value = 1737.6
$ws.Range("J80").Value = 2317.7778
$ws.Range("K80").Value = 1737.6
$ws.Range("L80").Value = 2317.7778
$ws.Range("M80").Value = -739.5999999999999
$ws.Range("N80").Value = -4313.7778
$ws.Range("H83").Value = 2012.421
$ws.Range("I83").Value = 1737.6
$ws.Range("J83").Value = 2317.7778
$ws.Range("K83").Value = 8688
$ws.Range("L83").Value = 11588.889
$ws.Range("M83").Value = -3696
$ws.Range("N83").Value = -21572.889
$ws.Range("H92").Value = 19474.75
$ws.Range("I92").Value = 20950
$ws.Range("J92").Value = 17999.5
$ws.Range("K92").Value = 20950
$ws.Range("L92").Value = 17999.5
$ws.Range("M92").Value = -19078
$ws.Range("N92").Value = -21743.5
$ws.Range("H102").Value = 1607.1904
$ws.Range("I102").Value = 631.9091
$ws.Range("K102").Value = 631.9091
$ws.Range("M102").Value = 990.0909
$ws.Range("H105").Value = 29999.5
$ws.Range("J105").Value = 29999.5
$ws.Range("L105").Value = 29999.5
$ws.Range("N105").Value = -36987.5
$ws.Range("H132").Value = 2720.4736
$ws.Range("I132").Value = 2023
$ws.Range("J132").Value = 3348.2
$ws.Range("K132").Value = 6069
$ws.Range("L132").Value = 10044.6
$ws.Range("M132").Value = -3539
$ws.Range("N132").Value = -15104.6

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 5266468.5
$ws.Range("J13").Value = 21950
$ws.Range("L13").Value = 21950
$ws.Range("N13").Value = -22230
$ws.Range("H22").Value = 14306.923
$ws.Range("I22").Value = 13498.75
$ws.Range("K22").Value = 13498.75
$ws.Range("M22").Value = -13203.75
$ws.Range("H26").Value = 10000
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5590
$ws.Range("H27").Value = 14306.923
$ws.Range("I27").Value = 13498.75
$ws.Range("K27").Value = 13498.75
$ws.Range("M27").Value = -13391.75
$ws.Range("H40").Value = 4995
$ws.Range("I40").Value = 4990
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 4990
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -4854
$ws.Range("N40").Value = -5272
$ws.Range("H46").Value = 2792.5557
$ws.Range("J46").Value = 3781.8333
$ws.Range("L46").Value = 3781.8333
$ws.Range("N46").Value = -4157.8333
$ws.Range("H60").Value = 53185
$ws.Range("J60").Value = 63870.5
$ws.Range("L60").Value = 63870.5
$ws.Range("N60").Value = -64888.5
$ws.Range("H68").Value = 1697.4
$ws.Range("I68").Value = 1496.75
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1496.75
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -747.75
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 1697.4
$ws.Range("I71").Value = 1496.75
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 7483.75
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -3739.75
$ws.Range("N71").Value = -19988
$ws.Range("H94").Value = 21975
$ws.Range("J94").Value = 21975
$ws.Range("L94").Value = 21975
$ws.Range("N94").Value = -23327
$ws.Range("H97").Value = 21521.5
$ws.Range("J97").Value = 21521.5
$ws.Range("L97").Value = 21521.5
$ws.Range("N97").Value = -23503.5
$ws.Range("H108").Value = 229500
$ws.Range("J108").Value = 229500
$ws.Range("L108").Value = 229500
$ws.Range("N108").Value = -237180
$ws.Range("H122").Value = 4478.1816
$ws.Range("I122").Value = 4084.5715
$ws.Range("K122").Value = 12253.7145
$ws.Range("M122").Value = -9803.7145
$ws.Range("H136").Value = 3244.1875
$ws.Range("I136").Value = 3063.5833
$ws.Range("J136").Value = 3786
$ws.Range("K136").Value = 9190.749899999999
$ws.Range("L136").Value = 11358
$ws.Range("M136").Value = -6640.749899999999
$ws.Range("N136").Value = -16458

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 47500
$ws.Range("J82").Value = 47500
$ws.Range("L82").Value = 47500
$ws.Range("N82").Value = -48266
$ws.Range("H85").Value = 47500
$ws.Range("J85").Value = 47500
$ws.Range("L85").Value = 47500
$ws.Range("N85").Value = -50152
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H97").Value = 57999
$ws.Range("J97").Value = 57999
$ws.Range("L97").Value = 57999
$ws.Range("N97").Value = -59981
$ws.Range("H126").Value = 1831.5883
$ws.Range("I126").Value = 1609.1333
$ws.Range("K126").Value = 4827.3999
$ws.Range("M126").Value = -2357.3999
$ws.Range("H132").Value = 42277
$ws.Range("I132").Value = 53496
$ws.Range("K132").Value = 160488
$ws.Range("M132").Value = -157958
$ws.Range("H136").Value = 1319.3636
$ws.Range("I136").Value = 1319.3636
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3958.0908
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1408.0908
$ws.Range("N136").Value = $null
